$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new journal entry as a new last row (row 20)
$ws.Cells.Item(20, 1).Value2 = 43929
$ws.Cells.Item(20, 2).Value = "Finalisation du projet, export, build et envoy de la 1.0 sur GitHub"

# Match formatting of the rows above (date format + bordered cell style)
$ws.Range("A13").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B13").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Update the active selection to B14, as captured in the saved view state
$ws.Range("B14").Select() | Out-Null
